$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.096772333333333
$ws.Range("H2").Value = 3.290317
$ws.Range("I2").Value = 0.2426185621302128
$ws.Range("J2").Value = 0.2426185621302128
$ws.Range("M2").Value = 1.819857
$ws.Range("N2").Value = 5.459571
$ws.Range("O2").Value = 0.01485317462584607
$ws.Range("P2").Value = 0.01485317462584607
$ws.Range("Q2").Value = 1.995968808223
$ws.Range("R2").Value = 17.963719274007
$ws.Range("S2").Value = 0.003603655870791735
$ws.Range("T2").Value = 0.003603655870791735

# Row 3
$ws.Range("G3").Value = 1.096772333333333
$ws.Range("H3").Value = 3.290317
$ws.Range("I3").Value = 0.2426185621302128
$ws.Range("J3").Value = 0.2426185621302128
$ws.Range("O3").Value = 0.726618572334523
$ws.Range("P3").Value = 0.7266185723345231
$ws.Range("Q3").Value = 97.64296471217332
$ws.Range("R3").Value = 878.7866824095599
$ws.Range("S3").Value = 0.17629115323691
$ws.Range("T3").Value = 0.17629115323691

# Row 4
$ws.Range("G4").Value = 1.096772333333333
$ws.Range("H4").Value = 3.290317
$ws.Range("I4").Value = 0.2426185621302128
$ws.Range("J4").Value = 0.2426185621302128
$ws.Range("M4").Value = 31.52924033333333
$ws.Range("N4").Value = 94.58772099999999
$ws.Range("O4").Value = 0.257333028084772
$ws.Range("P4").Value = 0.257333028084772
$ws.Range("Q4").Value = 34.58039848861744
$ws.Range("R4").Value = 311.223586397557
$ws.Range("S4").Value = 0.06243376926254107
$ws.Range("T4").Value = 0.06243376926254107

# Row 5
$ws.Range("G5").Value = 1.096772333333333
$ws.Range("H5").Value = 3.290317
$ws.Range("I5").Value = 0.2426185621302128
$ws.Range("J5").Value = 0.2426185621302128
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1464426666666667
$ws.Range("N5").Value = 0.439328
$ws.Range("O5").Value = 0.001195224954858853
$ws.Range("P5").Value = 0.001195224954858853
$ws.Range("Q5").Value = 0.1606142652195555
$ws.Range("R5").Value = 1.445528386976
$ws.Range("S5").Value = 0.0002899837599700033
$ws.Range("T5").Value = 0.0002899837599700034

# Row 6
$ws.Range("I6").Value = 0.03766810132102297
$ws.Range("J6").Value = 0.03766810132102297
$ws.Range("M6").Value = 1.819857
$ws.Range("N6").Value = 5.459571
$ws.Range("O6").Value = 0.01485317462584607
$ws.Range("P6").Value = 0.01485317462584607
$ws.Range("Q6").Value = 0.309887069817
$ws.Range("R6").Value = 2.788983628353
$ws.Range("S6").Value = 0.000559490886745217
$ws.Range("T6").Value = 0.0005594908867452171

# Row 7
$ws.Range("I7").Value = 0.03766810132102297
$ws.Range("J7").Value = 0.03766810132102297
$ws.Range("O7").Value = 0.726618572334523
$ws.Range("P7").Value = 0.7266185723345231
$ws.Range("S7").Value = 0.02737034200443387
$ws.Range("T7").Value = 0.02737034200443387

# Row 8
$ws.Range("I8").Value = 0.03766810132102297
$ws.Range("J8").Value = 0.03766810132102297
$ws.Range("M8").Value = 31.52924033333333
$ws.Range("N8").Value = 94.58772099999999
$ws.Range("O8").Value = 0.257333028084772
$ws.Range("P8").Value = 0.257333028084772
$ws.Range("Q8").Value = 5.368830573200333
$ws.Range("R8").Value = 48.319475158803
$ws.Range("S8").Value = 0.009693246575142842
$ws.Range("T8").Value = 0.009693246575142842

# Row 9
$ws.Range("I9").Value = 0.03766810132102297
$ws.Range("J9").Value = 0.03766810132102297
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1464426666666667
$ws.Range("N9").Value = 0.439328
$ws.Range("O9").Value = 0.001195224954858853
$ws.Range("P9").Value = 0.001195224954858853
$ws.Range("Q9").Value = 0.02493640372266667
$ws.Range("R9").Value = 0.224427633504
$ws.Range("S9").Value = 0.00004502185470103836
$ws.Range("T9").Value = 0.00004502185470103837

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5018676666666667
$ws.Range("H10").Value = 1.505603
$ws.Range("I10").Value = 0.1110188577571507
$ws.Range("J10").Value = 0.1110188577571507
$ws.Range("M10").Value = 1.819857
$ws.Range("N10").Value = 5.459571
$ws.Range("O10").Value = 0.01485317462584607
$ws.Range("P10").Value = 0.01485317462584607
$ws.Range("Q10").Value = 0.9133273862570002
$ws.Range("R10").Value = 8.219946476313
$ws.Range("S10").Value = 0.001648982481028925
$ws.Range("T10").Value = 0.001648982481028925

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5018676666666667
$ws.Range("H11").Value = 1.505603
$ws.Range("I11").Value = 0.1110188577571507
$ws.Range("J11").Value = 0.1110188577571507
$ws.Range("O11").Value = 0.726618572334523
$ws.Range("P11").Value = 0.7266185723345231
$ws.Range("Q11").Value = 44.68005380622667
$ws.Range("R11").Value = 402.12048425604
$ws.Range("S11").Value = 0.08066836392571035
$ws.Range("T11").Value = 0.08066836392571033

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5018676666666667
$ws.Range("H12").Value = 1.505603
$ws.Range("I12").Value = 0.1110188577571507
$ws.Range("J12").Value = 0.1110188577571507
$ws.Range("M12").Value = 31.52924033333333
$ws.Range("N12").Value = 94.58772099999999
$ws.Range("O12").Value = 0.257333028084772
$ws.Range("P12").Value = 0.257333028084772
$ws.Range("Q12").Value = 15.82350627786256
$ws.Range("R12").Value = 142.411556500763
$ws.Range("S12").Value = 0.02856881884116018
$ws.Range("T12").Value = 0.02856881884116018

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5018676666666667
$ws.Range("H13").Value = 1.505603
$ws.Range("I13").Value = 0.1110188577571507
$ws.Range("J13").Value = 0.1110188577571507
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1464426666666667
$ws.Range("N13").Value = 0.439328
$ws.Range("O13").Value = 0.001195224954858853
$ws.Range("P13").Value = 0.001195224954858853
$ws.Range("Q13").Value = 0.07349483942044445
$ws.Range("R13").Value = 0.661453554784
$ws.Range("S13").Value = 0.0001326925092512719
$ws.Range("T13").Value = 0.0001326925092512718

# Row 14
$ws.Range("G14").Value = 2.751641333333334
$ws.Range("H14").Value = 8.254924000000001
$ws.Range("I14").Value = 0.6086944787916135
$ws.Range("J14").Value = 0.6086944787916135
$ws.Range("M14").Value = 1.819857
$ws.Range("N14").Value = 5.459571
$ws.Range("O14").Value = 0.01485317462584607
$ws.Range("P14").Value = 0.01485317462584607
$ws.Range("Q14").Value = 5.007593741956001
$ws.Range("R14").Value = 45.06834367760401
$ws.Range("S14").Value = 0.00904104538728019
$ws.Range("T14").Value = 0.009041045387280192

# Row 15
$ws.Range("G15").Value = 2.751641333333334
$ws.Range("H15").Value = 8.254924000000001
$ws.Range("I15").Value = 0.6086944787916135
$ws.Range("J15").Value = 0.6086944787916135
$ws.Range("O15").Value = 0.726618572334523
$ws.Range("P15").Value = 0.7266185723345231
$ws.Range("Q15").Value = 244.9719139018133
$ws.Range("R15").Value = 2204.74722511632
$ws.Range("S15").Value = 0.4422887131674688
$ws.Range("T15").Value = 0.4422887131674689

# Row 16
$ws.Range("G16").Value = 2.751641333333334
$ws.Range("H16").Value = 8.254924000000001
$ws.Range("I16").Value = 0.6086944787916135
$ws.Range("J16").Value = 0.6086944787916135
$ws.Range("M16").Value = 31.52924033333333
$ws.Range("N16").Value = 94.58772099999999
$ws.Range("O16").Value = 0.257333028084772
$ws.Range("P16").Value = 0.257333028084772
$ws.Range("Q16").Value = 86.75716090980045
$ws.Range("R16").Value = 780.814448188204
$ws.Range("S16").Value = 0.156637193405928
$ws.Range("T16").Value = 0.156637193405928

# Row 17
$ws.Range("G17").Value = 2.751641333333334
$ws.Range("H17").Value = 8.254924000000001
$ws.Range("I17").Value = 0.6086944787916135
$ws.Range("J17").Value = 0.6086944787916135
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1464426666666667
$ws.Range("N17").Value = 0.439328
$ws.Range("O17").Value = 0.001195224954858853
$ws.Range("P17").Value = 0.001195224954858853
$ws.Range("Q17").Value = 0.4029576945635556
$ws.Range("R17").Value = 3.626619251072
$ws.Range("S17").Value = 0.000727526830936539
$ws.Range("T17").Value = 0.0007275268309365392
